# gateway_extraction_confidence_per_type analysis on TEST_DOCS
# Rebuild the sheet as "All Docs" (new sheetId) with an extra leading
# "Docs" column and six additional "Test Docs" rows.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheet (bumps the internal sheetId counter to 2)
#        and drop the old one, so "All Docs" ends up with sheetId="2".
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "All Docs"
$wb.Worksheets.Item("Tabelle1").Delete()

$ws = $wb.ActiveSheet

$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$xlVCenter = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$xlTop = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
$xlThin = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$xlContinuous = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$xlNone = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$xlLeft = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft
$xlRight = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight
$xlTopEdge = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop
$xlBottomEdge = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom
$xlAccent3 = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorAccent3

# --- 2. Header row -----------------------------------------------------
$ws.Cells.Item(1,1).Value2 = "Docs"
$ws.Cells.Item(1,2).Value2 = "Voting"
$ws.Cells.Item(1,3).Value2 = "Label"
$ws.Cells.Item(1,4).Value2 = "precision"
$ws.Cells.Item(1,5).Value2 = "recall"
$ws.Cells.Item(1,6).Value2 = "f1-score"
$ws.Cells.Item(1,7).Value2 = "confidence score"
$ws.Cells.Item(1,8).Value2 = "supports"

# --- 3. Data rows (2-7 = "All" docs, 8-13 = "Test Docs") --------------
# Columns: A Docs, B Voting, C Label, D precision, E recall, F f1-score,
#          G confidence score, H supports

# Row 2 - All / Full branch / XOR Gateway
$ws.Cells.Item(2,1).Value2 = "All"
$ws.Cells.Item(2,2).Value2 = "Full branch"
$ws.Cells.Item(2,3).Value2 = "XOR Gateway"
$ws.Cells.Item(2,4).Value2 = 0.92
$ws.Cells.Item(2,5).Value2 = 0.92
$ws.Cells.Item(2,6).Value2 = 0.92
$ws.Cells.Item(2,7).Value2 = 0.91
$ws.Cells.Item(2,8).Value2 = 66

# Row 3 - All / Full branch / AND Gateway
$ws.Cells.Item(3,1).Value2 = "All"
$ws.Cells.Item(3,2).Value2 = "Full branch"
$ws.Cells.Item(3,3).Value2 = "AND Gateway"
$ws.Cells.Item(3,4).Value2 = 0.75
$ws.Cells.Item(3,5).Value2 = 0.75
$ws.Cells.Item(3,6).Value2 = 0.75
$ws.Cells.Item(3,7).Value2 = 0.75
$ws.Cells.Item(3,8).Value2 = 8

# Row 4 - All / Full branch / Combined
$ws.Cells.Item(4,1).Value2 = "All"
$ws.Cells.Item(4,2).Value2 = "Full branch"
$ws.Cells.Item(4,3).Value2 = "Combined"
$ws.Cells.Item(4,4).Value2 = 0.9
$ws.Cells.Item(4,5).Value2 = 0.9
$ws.Cells.Item(4,6).Value2 = 0.9
$ws.Cells.Item(4,7).Formula = "=(G2*H2+G3*H3) /(H2+H3)"
$ws.Cells.Item(4,8).Value2 = 74

# Row 5 - All / Only start activity / XOR Gateway
$ws.Cells.Item(5,1).Value2 = "All"
$ws.Cells.Item(5,2).Value2 = "Only start activity"
$ws.Cells.Item(5,3).Value2 = "XOR Gateway"
$ws.Cells.Item(5,4).Value2 = 0.9
$ws.Cells.Item(5,5).Value2 = 0.91
$ws.Cells.Item(5,6).Value2 = 0.9
$ws.Cells.Item(5,7).Value2 = 0.98
$ws.Cells.Item(5,8).Value2 = 66

# Row 6 - All / Only start activity / AND Gateway
$ws.Cells.Item(6,1).Value2 = "All"
$ws.Cells.Item(6,2).Value2 = "Only start activity"
$ws.Cells.Item(6,3).Value2 = "AND Gateway"
$ws.Cells.Item(6,4).Value2 = 0.62
$ws.Cells.Item(6,5).Value2 = 0.62
$ws.Cells.Item(6,6).Value2 = 0.62
$ws.Cells.Item(6,7).Value2 = 1
$ws.Cells.Item(6,8).Value2 = 8

# Row 7 - All / Only start activity / Combined
$ws.Cells.Item(7,1).Value2 = "All"
$ws.Cells.Item(7,2).Value2 = "Only start activity"
$ws.Cells.Item(7,3).Value2 = "Combined"
$ws.Cells.Item(7,4).Value2 = 0.87
$ws.Cells.Item(7,5).Value2 = 0.88
$ws.Cells.Item(7,6).Value2 = 0.87
$ws.Cells.Item(7,7).Formula = "=(G5*H5+G6*H6) /(H5+H6)"
$ws.Cells.Item(7,8).Value2 = 74

# Row 8 - Test Docs / Full branch / XOR Gateway
$ws.Cells.Item(8,1).Value2 = "Test Docs"
$ws.Cells.Item(8,2).Value2 = "Full branch"
$ws.Cells.Item(8,3).Value2 = "XOR Gateway"
$ws.Cells.Item(8,4).Value2 = 1
$ws.Cells.Item(8,5).Value2 = 1
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 1
$ws.Cells.Item(8,8).Value2 = 12

# Row 9 - Test Docs / Full branch / AND Gateway
$ws.Cells.Item(9,1).Value2 = "Test Docs"
$ws.Cells.Item(9,2).Value2 = "Full branch"
$ws.Cells.Item(9,3).Value2 = "AND Gateway"
$ws.Cells.Item(9,4).Value2 = 1
$ws.Cells.Item(9,5).Value2 = 1
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 1
$ws.Cells.Item(9,8).Value2 = 1

# Row 10 - Test Docs / Full branch / Combined
$ws.Cells.Item(10,1).Value2 = "Test Docs"
$ws.Cells.Item(10,2).Value2 = "Full branch"
$ws.Cells.Item(10,3).Value2 = "Combined"
$ws.Cells.Item(10,4).Value2 = 1
$ws.Cells.Item(10,5).Value2 = 1
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 1
$ws.Cells.Item(10,8).Value2 = 13

# Row 11 - Test Docs / Only start activity / XOR Gateway
$ws.Cells.Item(11,1).Value2 = "Test Docs"
$ws.Cells.Item(11,2).Value2 = "Only start activity"
$ws.Cells.Item(11,3).Value2 = "XOR Gateway"
$ws.Cells.Item(11,4).Value2 = 1
$ws.Cells.Item(11,5).Value2 = 1
$ws.Cells.Item(11,6).Value2 = 1
$ws.Cells.Item(11,7).Value2 = 1
$ws.Cells.Item(11,8).Value2 = 12

# Row 12 - Test Docs / Only start activity / AND Gateway
$ws.Cells.Item(12,1).Value2 = "Test Docs"
$ws.Cells.Item(12,2).Value2 = "Only start activity"
$ws.Cells.Item(12,3).Value2 = "AND Gateway"
$ws.Cells.Item(12,4).Value2 = 1
$ws.Cells.Item(12,5).Value2 = 1
$ws.Cells.Item(12,6).Value2 = 1
$ws.Cells.Item(12,7).Value2 = 1
$ws.Cells.Item(12,8).Value2 = 1

# Row 13 - Test Docs / Only start activity / Combined
$ws.Cells.Item(13,1).Value2 = "Test Docs"
$ws.Cells.Item(13,2).Value2 = "Only start activity"
$ws.Cells.Item(13,3).Value2 = "Combined"
$ws.Cells.Item(13,4).Value2 = 1
$ws.Cells.Item(13,5).Value2 = 1
$ws.Cells.Item(13,6).Value2 = 1
$ws.Cells.Item(13,7).Value2 = 1
$ws.Cells.Item(13,8).Value2 = 13

# --- 4. Column widths ---------------------------------------------------
$ws.Range("B1").ColumnWidth = 17.88671875
$ws.Range("C1").ColumnWidth = 12.77734375
$ws.Range("D1").ColumnWidth = 8.88671875
$ws.Range("E1").ColumnWidth = 5.6640625
$ws.Range("F1").ColumnWidth = 7.88671875
$ws.Range("G1").ColumnWidth = 15.6640625
$ws.Range("H1").ColumnWidth = 8.44140625

# --- 5. Formatting -------------------------------------------------------

# Whole table: centered horizontally
$ws.Range("A1:H13").HorizontalAlignment = $xlCenter

# Header row: bold, vertical top
$ws.Range("A1:H1").Font.Bold = $true
$ws.Range("A1:H1").VerticalAlignment = $xlTop

# A1:B1 (Docs / Voting headers) - full thin box border (reuse original header style)
$ws.Range("A1:B1").Borders.LineStyle = $xlContinuous
$ws.Range("A1:B1").Borders.Weight = $xlThin

# C1:F1 (Label / precision / recall / f1-score headers) - left+right+top border
$rng = $ws.Range("C1:F1")
$rng.Borders.Item($xlLeft).LineStyle = $xlContinuous
$rng.Borders.Item($xlLeft).Weight = $xlThin
$rng.Borders.Item($xlRight).LineStyle = $xlContinuous
$rng.Borders.Item($xlRight).Weight = $xlThin
$rng.Borders.Item($xlTopEdge).LineStyle = $xlContinuous
$rng.Borders.Item($xlTopEdge).Weight = $xlThin

# G1 (confidence score header) same border treatment, accent color + bold
$rng = $ws.Range("G1")
$rng.Borders.Item($xlLeft).LineStyle = $xlContinuous
$rng.Borders.Item($xlLeft).Weight = $xlThin
$rng.Borders.Item($xlRight).LineStyle = $xlContinuous
$rng.Borders.Item($xlRight).Weight = $xlThin
$rng.Borders.Item($xlTopEdge).LineStyle = $xlContinuous
$rng.Borders.Item($xlTopEdge).Weight = $xlThin
$rng.Font.ThemeColor = $xlAccent3
$rng.Font.Bold = $true

# H1 (supports header) full thin box border like original
$ws.Range("H1").Borders.LineStyle = $xlContinuous
$ws.Range("H1").Borders.Weight = $xlThin

# --- Column B/C (Voting / Label) - left border all data rows ----------
$ws.Range("B2:C13").Borders.Item($xlLeft).LineStyle = $xlContinuous
$ws.Range("B2:C13").Borders.Item($xlLeft).Weight = $xlThin
# group-top rows (2,5,8,11) also get a top border
foreach ($r in 2,5,8,11) {
    $rr = $ws.Range("B" + $r + ":C" + $r)
    $rr.Borders.Item($xlTopEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlTopEdge).Weight = $xlThin
}
# group-bottom rows (4,7,10,13) also get a bottom border
foreach ($r in 4,7,10,13) {
    $rr = $ws.Range("B" + $r + ":C" + $r)
    $rr.Borders.Item($xlBottomEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlBottomEdge).Weight = $xlThin
}

# --- Columns D:F (precision/recall/f1) - left+right border, vertical-center
$ws.Range("D2:F13").VerticalAlignment = $xlVCenter
$ws.Range("D2:F13").Borders.Item($xlLeft).LineStyle = $xlContinuous
$ws.Range("D2:F13").Borders.Item($xlLeft).Weight = $xlThin
$ws.Range("D2:F13").Borders.Item($xlRight).LineStyle = $xlContinuous
$ws.Range("D2:F13").Borders.Item($xlRight).Weight = $xlThin
foreach ($r in 2,5,8,11) {
    $rr = $ws.Range("D" + $r + ":F" + $r)
    $rr.Borders.Item($xlTopEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlTopEdge).Weight = $xlThin
}
foreach ($r in 4,7,10,13) {
    $rr = $ws.Range("D" + $r + ":F" + $r)
    $rr.Borders.Item($xlBottomEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlBottomEdge).Weight = $xlThin
}

# --- Column G (confidence score) - right border, vertical-center, accent color on non-formula rows
$ws.Range("G2:G13").VerticalAlignment = $xlVCenter
$ws.Range("G2:G13").Borders.Item($xlRight).LineStyle = $xlContinuous
$ws.Range("G2:G13").Borders.Item($xlRight).Weight = $xlThin
foreach ($r in 2,5,8,11) {
    $rr = $ws.Range("G" + $r)
    $rr.Borders.Item($xlTopEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlTopEdge).Weight = $xlThin
}
foreach ($r in 4,7,10,13) {
    $rr = $ws.Range("G" + $r)
    $rr.Borders.Item($xlBottomEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlBottomEdge).Weight = $xlThin
}
# the "per-row" confidence score values (not the combined-formula rows) get the accent font color
$ws.Range("G2").Font.ThemeColor = $xlAccent3
$ws.Range("G3").Font.ThemeColor = $xlAccent3
$ws.Range("G5").Font.ThemeColor = $xlAccent3
$ws.Range("G6").Font.ThemeColor = $xlAccent3
$ws.Range("G7").Font.ThemeColor = $xlAccent3

# --- Column H (supports) - left+right border, vertical-center
$ws.Range("H2:H13").VerticalAlignment = $xlVCenter
$ws.Range("H2:H13").Borders.Item($xlLeft).LineStyle = $xlContinuous
$ws.Range("H2:H13").Borders.Item($xlLeft).Weight = $xlThin
$ws.Range("H2:H13").Borders.Item($xlRight).LineStyle = $xlContinuous
$ws.Range("H2:H13").Borders.Item($xlRight).Weight = $xlThin
foreach ($r in 2,5,8,11) {
    $rr = $ws.Range("H" + $r)
    $rr.Borders.Item($xlTopEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlTopEdge).Weight = $xlThin
}
foreach ($r in 4,7,10,13) {
    $rr = $ws.Range("H" + $r)
    $rr.Borders.Item($xlBottomEdge).LineStyle = $xlContinuous
    $rr.Borders.Item($xlBottomEdge).Weight = $xlThin
}

# --- Column A (Docs) - only a right border, shaded fill for the "Test Docs" block
$ws.Range("A2:A13").Borders.Item($xlRight).LineStyle = $xlContinuous
$ws.Range("A2:A13").Borders.Item($xlRight).Weight = $xlThin
$ws.Range("A7").Borders.Item($xlBottomEdge).LineStyle = $xlContinuous
$ws.Range("A7").Borders.Item($xlBottomEdge).Weight = $xlThin
$ws.Range("A13").Borders.Item($xlBottomEdge).LineStyle = $xlContinuous
$ws.Range("A13").Borders.Item($xlBottomEdge).Weight = $xlThin

# --- 6. Comment: move from (old) F1 to G1 (confidence score header) ----
$commentText = "Computed with gateway_extraction_confidence_per_type"
$ws.Range("G1").AddCommentThreaded($commentText)

# --- 7. Workbook window view size ---------------------------------------
$excel.Width = 14400
$excel.Height = 15600
$excel.Left = -28800
$excel.Top = 0

$ws.Cells.Item(1,1).Select()
